$wb = $excel.ActiveWorkbook

# Sheet "Overview": File Name | zh-cn (Status) | de-de (Status) | Latest Handoff Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-03-24 22:38:05"

# Sheet "zh-cn": ... | Status (C) | ... | Latest Handoff Datetime (E) | ...
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("E2").Value = "2016-03-24 22:38:00"

# Sheet "de-de": ... | Status (C) | ... | Latest Handoff Datetime (E) | ...
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("E2").Value = "2016-03-24 22:38:05"
